$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that changed from 45189 (2023-09-20)
# to 45190 (2023-09-21) for every data row (rows 2 through 27).
$ws.Range("C2:C27").Value = 45190
